$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.219.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.784.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.47%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3780"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07498"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.17%  "

$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.472"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.786.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.102"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06671"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("E19").Value = "  -3.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.613"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.219.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.417"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.547"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.989.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.93"
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = "  -2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.097"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08697"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.663"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6955"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("E38").Value = "  -4.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2202"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.54%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.814"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06337"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02339"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.242"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6514"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.848"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.147"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07130"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "

